$wb = $excel.ActiveWorkbook

# --- Sheet 1 ---
$ws = $wb.Worksheets.Item(1)
$ws.Name = "summ27578491"
$ws.Rows.Item(13).Delete()
$ws.Range("B2").Value = [double]"7246.894834399931"
$ws.Range("C2").Value = [double]"6.668222704437674e-05"
$ws.Range("B3").Value = [double]"-1511.056910150845"
$ws.Range("C3").Value = [double]"0.2368628427982399"
$ws.Range("B4").Value = [double]"-887.9949243176047"
$ws.Range("C4").Value = [double]"0.4886611442110637"
$ws.Range("B5").Value = [double]"-433.1137178985263"
$ws.Range("C5").Value = [double]"0.1392699277534862"
$ws.Range("B6").Value = [double]"67.99289478421247"
$ws.Range("C6").Value = [double]"0.5868797052532284"
$ws.Range("B7").Value = [double]"-776.9412922652144"
$ws.Range("C7").Value = [double]"0.008076650234334256"
$ws.Range("B8").Value = [double]"-38.72583806043021"
$ws.Range("C8").Value = [double]"0.003399738128581711"
$ws.Range("B9").Value = [double]"65.6594456547941"
$ws.Range("C9").Value = [double]"0.5955090559157934"
$ws.Range("B10").Value = [double]"610.6971999422312"
$ws.Range("C10").Value = [double]"1.736559383233683e-15"
$ws.Range("B11").Value = [double]"-0.1589794990991311"
$ws.Range("C11").Value = [double]"0.2203628878800351"
$ws.Range("B12").Value = [double]"0.0001414622857089593"
$ws.Range("C12").Value = [double]"0.3306007996743626"
$ws.Range("B13").Value = [double]"18.4404569750178"
$ws.Range("C13").Value = [double]"0.001767917100464464"
$ws.Range("B14").Value = [double]"-5755.993188725779"
$ws.Range("C14").Value = [double]"0.004398736140770852"
$ws.Range("B15").Value = [double]"-2558.554697915431"
$ws.Range("C15").Value = [double]"0.03186318411949117"

# --- Sheet 2 ---
$ws = $wb.Worksheets.Item(2)
$ws.Name = "summ27789082"
$ws.Rows.Item(13).Delete()
$ws.Range("B2").Value = [double]"5347.853526777185"
$ws.Range("C2").Value = [double]"0.002575306780658934"
$ws.Range("B3").Value = [double]"-289.21607960489"
$ws.Range("C3").Value = [double]"0.8117233609095895"
$ws.Range("B4").Value = [double]"468.148739808108"
$ws.Range("C4").Value = [double]"0.7012174025903058"
$ws.Range("B5").Value = [double]"-178.228386775684"
$ws.Range("C5").Value = [double]"0.5506645844286107"
$ws.Range("B6").Value = [double]"36.75058121570274"
$ws.Range("C6").Value = [double]"0.7728724880354736"
$ws.Range("B7").Value = [double]"-444.0625559443728"
$ws.Range("C7").Value = [double]"0.1366779886661871"
$ws.Range("B8").Value = [double]"-45.18061253613169"
$ws.Range("C8").Value = [double]"0.0008959452405628962"
$ws.Range("B9").Value = [double]"10.09344597664979"
$ws.Range("C9").Value = [double]"0.9347273514862932"
$ws.Range("B10").Value = [double]"659.623071399121"
$ws.Range("C10").Value = [double]"9.366726620617052e-19"
$ws.Range("B11").Value = [double]"-0.14906580618655"
$ws.Range("C11").Value = [double]"0.2726013434047908"
$ws.Range("B12").Value = [double]"9.757404960065717e-05"
$ws.Range("C12").Value = [double]"0.5261085245462809"
$ws.Range("B13").Value = [double]"23.58714551416669"
$ws.Range("C13").Value = [double]"0.00017291476986433"
$ws.Range("B14").Value = [double]"-4976.778418263709"
$ws.Range("C14").Value = [double]"0.02033376970769971"
$ws.Range("B15").Value = [double]"-2111.495404910145"
$ws.Range("C15").Value = [double]"0.08146086969357827"

# --- Sheet 3 ---
$ws = $wb.Worksheets.Item(3)
$ws.Name = "summ28055802"
$ws.Rows.Item(13).Delete()
$ws.Range("B2").Value = [double]"4570.726447510337"
$ws.Range("C2").Value = [double]"0.00657676676495771"
$ws.Range("B3").Value = [double]"-133.5472730630868"
$ws.Range("C3").Value = [double]"0.9057753772789547"
$ws.Range("B4").Value = [double]"459.898871062388"
$ws.Range("C4").Value = [double]"0.6849045342205746"
$ws.Range("B5").Value = [double]"-585.1000945753875"
$ws.Range("C5").Value = [double]"0.04570094370911403"
$ws.Range("B6").Value = [double]"-172.9241291034691"
$ws.Range("C6").Value = [double]"0.1725513584601903"
$ws.Range("B7").Value = [double]"-332.989494904044"
$ws.Range("C7").Value = [double]"0.2532616934908626"
$ws.Range("B8").Value = [double]"-39.74550243311096"
$ws.Range("C8").Value = [double]"0.002410367779414926"
$ws.Range("B9").Value = [double]"-91.30101089480061"
$ws.Range("C9").Value = [double]"0.4482628230323336"
$ws.Range("B10").Value = [double]"739.4975994525405"
$ws.Range("C10").Value = [double]"2.069787090572408e-23"
$ws.Range("B11").Value = [double]"-0.04527971707977282"
$ws.Range("C11").Value = [double]"0.7240729768980605"
$ws.Range("B12").Value = [double]"9.588181929198386e-05"
$ws.Range("C12").Value = [double]"0.5030869094959336"
$ws.Range("B13").Value = [double]"30.58925221459952"
$ws.Range("C13").Value = [double]"2.881305033068816e-07"
$ws.Range("B14").Value = [double]"-5451.831087059209"
$ws.Range("C14").Value = [double]"0.00611172533334087"
$ws.Range("B15").Value = [double]"-1966.563444089793"
$ws.Range("C15").Value = [double]"0.09576829702806135"

# --- Sheet 4 ---
$ws = $wb.Worksheets.Item(4)
$ws.Name = "summ28325781"
$ws.Rows.Item(13).Delete()
$ws.Range("B2").Value = [double]"5777.880232228314"
$ws.Range("C2").Value = [double]"0.0005621194674598266"
$ws.Range("B3").Value = [double]"-331.9019303789574"
$ws.Range("C3").Value = [double]"0.7690101439157259"
$ws.Range("B4").Value = [double]"170.4500755531752"
$ws.Range("C4").Value = [double]"0.8806029386230403"
$ws.Range("B5").Value = [double]"-360.8786659658905"
$ws.Range("C5").Value = [double]"0.2162734564810117"
$ws.Range("B6").Value = [double]"-21.04494373091909"
$ws.Range("C6").Value = [double]"0.8632575003931313"
$ws.Range("B7").Value = [double]"-362.4399170932361"
$ws.Range("C7").Value = [double]"0.2131962517584552"
$ws.Range("B8").Value = [double]"-31.88695278430173"
$ws.Range("C8").Value = [double]"0.01469081539137169"
$ws.Range("B9").Value = [double]"33.06578880675528"
$ws.Range("C9").Value = [double]"0.7873226091689578"
$ws.Range("B10").Value = [double]"691.4177832118119"
$ws.Range("C10").Value = [double]"4.22656960092371e-21"
$ws.Range("B11").Value = [double]"-0.1592376251370658"
$ws.Range("C11").Value = [double]"0.2160157493607978"
$ws.Range("B12").Value = [double]"0.0001995173830019937"
$ws.Range("C12").Value = [double]"0.1716888952344015"
$ws.Range("B13").Value = [double]"14.94559100506491"
$ws.Range("C13").Value = [double]"0.01470971635457912"
$ws.Range("B14").Value = [double]"-6809.801249694048"
$ws.Range("C14").Value = [double]"0.001155098088651693"
$ws.Range("B15").Value = [double]"-2466.201828664252"
$ws.Range("C15").Value = [double]"0.03539899322496337"

# --- Sheet 5 ---
$ws = $wb.Worksheets.Item(5)
$ws.Name = "summ28588249"
$ws.Rows.Item(13).Delete()
$ws.Range("B2").Value = [double]"6759.606483079187"
$ws.Range("C2").Value = [double]"8.953823533017366e-05"
$ws.Range("B3").Value = [double]"-1576.533125521908"
$ws.Range("C3").Value = [double]"0.2026688850466428"
$ws.Range("B4").Value = [double]"-719.447526970122"
$ws.Range("C4").Value = [double]"0.5623004783214336"
$ws.Range("B5").Value = [double]"-462.8980275640846"
$ws.Range("C5").Value = [double]"0.114639985239854"
$ws.Range("B6").Value = [double]"52.60213126442421"
$ws.Range("C6").Value = [double]"0.667519371280094"
$ws.Range("B7").Value = [double]"-404.2848408716868"
$ws.Range("C7").Value = [double]"0.167490821772482"
$ws.Range("B8").Value = [double]"-48.25649637041712"
$ws.Range("C8").Value = [double]"0.0002408219886606014"
$ws.Range("B9").Value = [double]"-66.23225428913048"
$ws.Range("C9").Value = [double]"0.5902712352224828"
$ws.Range("B10").Value = [double]"691.8318356417321"
$ws.Range("C10").Value = [double]"3.61869153537552e-21"
$ws.Range("B11").Value = [double]"-0.1150476492042373"
$ws.Range("C11").Value = [double]"0.3779891651288306"
$ws.Range("B12").Value = [double]"9.960888539576662e-05"
$ws.Range("C12").Value = [double]"0.5001874865994511"
$ws.Range("B13").Value = [double]"23.86812965101461"
$ws.Range("C13").Value = [double]"3.58372225339273e-05"
$ws.Range("B14").Value = [double]"-5556.305306707854"
$ws.Range("C14").Value = [double]"0.007453858159417466"
$ws.Range("B15").Value = [double]"-2416.09576276948"
$ws.Range("C15").Value = [double]"0.04677921739250908"

# --- Sheet 6 ---
$ws = $wb.Worksheets.Item(6)
$ws.Name = "summ28852960"
$ws.Rows.Item(13).Delete()
$ws.Range("B2").Value = [double]"5078.428157483289"
$ws.Range("C2").Value = [double]"0.008788065326140192"
$ws.Range("B3").Value = [double]"830.9094143629513"
$ws.Range("C3").Value = [double]"0.5781531625619157"
$ws.Range("B4").Value = [double]"1443.193429794025"
$ws.Range("C4").Value = [double]"0.3349106258088419"
$ws.Range("B5").Value = [double]"-563.4860740219445"
$ws.Range("C5").Value = [double]"0.05581105582682931"
$ws.Range("B6").Value = [double]"-2.414938342266225"
$ws.Range("C6").Value = [double]"0.9842161432362013"
$ws.Range("B7").Value = [double]"-688.1408207280308"
$ws.Range("C7").Value = [double]"0.01852997406510874"
$ws.Range("B8").Value = [double]"-53.62397103904085"
$ws.Range("C8").Value = [double]"4.390272714037106e-05"
$ws.Range("B9").Value = [double]"-27.38054638544143"
$ws.Range("C9").Value = [double]"0.8220581198107775"
$ws.Range("B10").Value = [double]"664.2043740377558"
$ws.Range("C10").Value = [double]"1.250542792493233e-18"
$ws.Range("B11").Value = [double]"-0.10441735050713"
$ws.Range("C11").Value = [double]"0.4240486556944376"
$ws.Range("B12").Value = [double]"8.092953770983942e-05"
$ws.Range("C12").Value = [double]"0.5817422715330062"
$ws.Range("B13").Value = [double]"23.06046215062793"
$ws.Range("C13").Value = [double]"6.481946134572963e-05"
$ws.Range("B14").Value = [double]"-6186.019853040545"
$ws.Range("C14").Value = [double]"0.00259356581176264"
$ws.Range("B15").Value = [double]"-1818.52859888708"
$ws.Range("C15").Value = [double]"0.1308084987487612"

# --- Sheet 7 ---
$ws = $wb.Worksheets.Item(7)
$ws.Name = "summ29109108"
$ws.Rows.Item(13).Delete()
$ws.Range("B2").Value = [double]"6739.265602263369"
$ws.Range("C2").Value = [double]"9.152518033534386e-05"
$ws.Range("B3").Value = [double]"-231.4165300149948"
$ws.Range("C3").Value = [double]"0.8401983770591083"
$ws.Range("B4").Value = [double]"124.5343936359698"
$ws.Range("C4").Value = [double]"0.9139930742419271"
$ws.Range("B5").Value = [double]"-250.2463287060955"
$ws.Range("C5").Value = [double]"0.3994152552757492"
$ws.Range("B6").Value = [double]"-4.072471263793034"
$ws.Range("C6").Value = [double]"0.9746057566015531"
$ws.Range("B7").Value = [double]"-389.9139243588982"
$ws.Range("C7").Value = [double]"0.1880677092151543"
$ws.Range("B8").Value = [double]"-45.55516945734027"
$ws.Range("C8").Value = [double]"0.0007003769650646962"
$ws.Range("B9").Value = [double]"49.71670071574482"
$ws.Range("C9").Value = [double]"0.6868585737058298"
$ws.Range("B10").Value = [double]"586.8293827090729"
$ws.Range("C10").Value = [double]"2.144310280516263e-15"
$ws.Range("B11").Value = [double]"-0.1281641205409248"
$ws.Range("C11").Value = [double]"0.3220380460527162"
$ws.Range("B12").Value = [double]"7.694726499871585e-05"
$ws.Range("C12").Value = [double]"0.6013313134185543"
$ws.Range("B13").Value = [double]"16.91098336580277"
$ws.Range("C13").Value = [double]"0.0042742255095024"
$ws.Range("B14").Value = [double]"-6017.664395527881"
$ws.Range("C14").Value = [double]"0.005169382130807597"
$ws.Range("B15").Value = [double]"-2343.411843496253"
$ws.Range("C15").Value = [double]"0.05274576982798843"

# --- Sheet 8 ---
$ws = $wb.Worksheets.Item(8)
$ws.Name = "summ29364017"
$ws.Rows.Item(13).Delete()
$ws.Range("B2").Value = [double]"4130.053876226246"
$ws.Range("C2").Value = [double]"0.01594020661993227"
$ws.Range("B3").Value = [double]"875.4442394023671"
$ws.Range("C3").Value = [double]"0.4580722896739772"
$ws.Range("B4").Value = [double]"1479.46800893869"
$ws.Range("C4").Value = [double]"0.2117250147915611"
$ws.Range("B5").Value = [double]"-212.3768171316408"
$ws.Range("C5").Value = [double]"0.462242838402113"
$ws.Range("B6").Value = [double]"-59.13823196969238"
$ws.Range("C6").Value = [double]"0.6293313719741243"
$ws.Range("B7").Value = [double]"-498.4333204061813"
$ws.Range("C7").Value = [double]"0.08475876006266919"
$ws.Range("B8").Value = [double]"-49.04610787090721"
$ws.Range("C8").Value = [double]"0.0002061099195876843"
$ws.Range("B9").Value = [double]"61.87710762124166"
$ws.Range("C9").Value = [double]"0.604589332987207"
$ws.Range("B10").Value = [double]"715.1870916617509"
$ws.Range("C10").Value = [double]"3.836093402639869e-23"
$ws.Range("B11").Value = [double]"-0.1663312001822424"
$ws.Range("C11").Value = [double]"0.1938313590456818"
$ws.Range("B12").Value = [double]"0.0002339428796325996"
$ws.Range("C12").Value = [double]"0.1034724933899865"
$ws.Range("B13").Value = [double]"19.94304693757027"
$ws.Range("C13").Value = [double]"0.00049268828864712"
$ws.Range("B14").Value = [double]"-6751.77650701825"
$ws.Range("C14").Value = [double]"0.0008023886696772479"
$ws.Range("B15").Value = [double]"-1398.848009922712"
$ws.Range("C15").Value = [double]"0.2404082704891742"

# --- Sheet 9 ---
$ws = $wb.Worksheets.Item(9)
$ws.Name = "summ29639373"
$ws.Rows.Item(13).Delete()
$ws.Range("B2").Value = [double]"4502.614873481172"
$ws.Range("C2").Value = [double]"0.005897968480136383"
$ws.Range("B3").Value = [double]"-467.6675877417566"
$ws.Range("C3").Value = [double]"0.6718435345794725"
$ws.Range("B4").Value = [double]"214.416198370282"
$ws.Range("C4").Value = [double]"0.8466700200072367"
$ws.Range("B5").Value = [double]"115.4516593853699"
$ws.Range("C5").Value = [double]"0.6840474660575628"
$ws.Range("B6").Value = [double]"57.34581107968653"
$ws.Range("C6").Value = [double]"0.6343859468384156"
$ws.Range("B7").Value = [double]"-869.0723861763797"
$ws.Range("C7").Value = [double]"0.002329868133970071"
$ws.Range("B8").Value = [double]"-31.37954355357687"
$ws.Range("C8").Value = [double]"0.0145501873331041"
$ws.Range("B9").Value = [double]"18.65488186101018"
$ws.Range("C9").Value = [double]"0.8741223747193668"
$ws.Range("B10").Value = [double]"670.1436113736358"
$ws.Range("C10").Value = [double]"1.464626334721252e-20"
$ws.Range("B11").Value = [double]"-0.1163539426031062"
$ws.Range("C11").Value = [double]"0.3558834770096129"
$ws.Range("B12").Value = [double]"0.0001419877783299893"
$ws.Range("C12").Value = [double]"0.3164773977927807"
$ws.Range("B13").Value = [double]"20.16276462069916"
$ws.Range("C13").Value = [double]"0.0003279673051605143"
$ws.Range("B14").Value = [double]"-5677.254442679677"
$ws.Range("C14").Value = [double]"0.003855339877710808"
$ws.Range("B15").Value = [double]"-1160.083440235122"
$ws.Range("C15").Value = [double]"0.3128582731885057"

